$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 22:32"

# Update country rows whose data changed between the two data pulls.
# Row layout: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#             E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 6994855
$ws.Range("C4").Value = 27452
$ws.Range("D4").Value = 4243371
$ws.Range("E4").Value = 2547435
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 225
$ws.Range("H4").Value = 204049

# Row 5: India
$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 5485612
$ws.Range("C5").Value = 87382
$ws.Range("D5").Value = 4392650
$ws.Range("E5").Value = 1005053
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1135
$ws.Range("H5").Value = 87909

# Row 25: Alemania
$ws.Range("A25").Value = "Alemania"
$ws.Range("B25").Value = 273477
$ws.Range("C25").Value = 1169
$ws.Range("D25").Value = 243500
$ws.Range("E25").Value = 20507
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 9470

# Row 27: Israel
$ws.Range("A27").Value = "Israel"
$ws.Range("B27").Value = 187902
$ws.Range("C27").Value = 4300
$ws.Range("D27").Value = 134069
$ws.Range("E27").Value = 52577
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 30
$ws.Range("H27").Value = 1256

# Row 49: Bielorrusia
$ws.Range("A49").Value = "Bielorrusia"
$ws.Range("B49").Value = 75674
$ws.Range("C49").Value = 213
$ws.Range("D49").Value = 73265
$ws.Range("E49").Value = 1629
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 780

# Row 83: Costa de Marfil
$ws.Range("A83").Value = "Costa de Marfil"
$ws.Range("B83").Value = 19320
$ws.Range("C83").Value = 51
$ws.Range("D83").Value = 18460
$ws.Range("E83").Value = 740
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 120

# Row 126: Ruanda
$ws.Range("A126").Value = "Ruanda"
$ws.Range("B126").Value = 4711
$ws.Range("C126").Value = 22
$ws.Range("D126").Value = 2961
$ws.Range("E126").Value = 1724
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 26

# Row 127: Surinam
$ws.Range("A127").Value = "Surinam"
$ws.Range("B127").Value = 4709
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 4383
$ws.Range("E127").Value = 229
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 97

# Row 204: Santa Lucia
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("B204").Value = 27
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 26
$ws.Range("E204").Value = 1
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

# Row 205: Timor Oriental
$ws.Range("A205").Value = "Timor Oriental"
$ws.Range("B205").Value = 27
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 26
$ws.Range("E205").Value = 1
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# Row 214: Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

# Row 215: Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
